# Fruta / hortaliza, semanal
# Insert a new week's data row at row 7, push the old row 7 -> row 8 and old
# row 8 -> row 9, updating the dates/values to match the new weekly snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 gets what used to be row 8's data (unchanged values).
$ws.Range("A9").Value = 3
$ws.Range("B9").Value = "Femacal de La Calera"
$ws.Range("C9").Value = "Coquimbo"
$ws.Range("D9").Value = 44209
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100101
$ws.Range("H9").Value = "Berries"
$ws.Range("I9").Value = 100101004
$ws.Range("J9").Value = "Frambuesa"
$ws.Range("K9").Value = "Sin especificar"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 58
$ws.Range("N9").Value = 6000
$ws.Range("O9").Value = 6000
$ws.Range("P9").Value = 6000
$ws.Range("Q9").Value = '$/bandeja 2 kilos'
$ws.Range("R9").Value = "Provincia de Curicó"
$ws.Range("S9").Value = 3000
$ws.Range("T9").Value = 2

# copy the date number format from D8 to the new D9 cell
$ws.Range("D8").Copy()
$ws.Range("D9").PasteSpecial(-4122)  # xlPasteFormats

# Row 8 now becomes what used to be row 7's values, but with an updated date.
$ws.Range("D8").Value = 44585
$ws.Range("M8").Value = 160
$ws.Range("N8").Value = 6500
$ws.Range("O8").Value = 7000
$ws.Range("P8").Value = 6750
$ws.Range("S8").Value = 3375
$ws.Range("T8").Value = 2

# Row 7 keeps the same values, just a newer date for the latest week.
$ws.Range("D7").Value = 44588
